$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Title string in B1 on both sheets lost the space before "Repayment".
$wsInput.Range("B1").Value  = "1015-MS-EI-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-LateRepayment"
$wsOutput.Range("B1").Value = "1015-MS-EI-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-LateRepayment"

# Active tab moves from the second sheet ("ProductLoanOutput") to the first
# ("ProductLoanInput"), and both sheets end up with B1 selected/topmost
# instead of their previous scroll position/selection. Select the
# soon-to-be-inactive sheet first so the final selection below wins as the
# active tab.
[void]$wsOutput.Range("B1").Select()
[void]$wsInput.Range("B1").Select()
